$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Cells whose text looks numeric (e.g. NPI numbers) need NumberFormat forced to
# Text before assignment so Excel stores them as shared strings (matching the
# source data) instead of auto-converting to a Number type. ClearFormats()
# afterwards drops the temporary format so the cell keeps the workbook default
# style (no "s" attribute), matching the original cells around it.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-TextValue $ws.Range('I42') '1134117690'
Set-TextValue $ws.Range('I43') '1134117690'
Set-TextValue $ws.Range('I51') '1952386393'
Set-TextValue $ws.Range('I57') '1134377922'
Set-TextValue $ws.Range('I58') '1134377922'
Set-TextValue $ws.Range('I106') '1578755476'
Set-TextValue $ws.Range('I107') '1578755476'
Set-TextValue $ws.Range('I112') '1568426302'
Set-TextValue $ws.Range('I120') '1902847445'
Set-TextValue $ws.Range('I121') '1902847445'
Set-TextValue $ws.Range('I126') '1619918042'
Set-TextValue $ws.Range('I127') '1619918042'
Set-TextValue $ws.Range('I134') '1386930493'
Set-TextValue $ws.Range('I135') '1386930493'
Set-TextValue $ws.Range('I145') '1275828089'
Set-TextValue $ws.Range('I146') '1275828089'
Set-TextValue $ws.Range('I157') '1184650145'
Set-TextValue $ws.Range('I158') '1184650145'
Set-TextValue $ws.Range('I161') '1528286606'
Set-TextValue $ws.Range('I162') '1528286606'
Set-TextValue $ws.Range('I171') '1740571116'
Set-TextValue $ws.Range('I172') '1740571116'
Set-TextValue $ws.Range('I175') '1588807416'
Set-TextValue $ws.Range('I178') '1689641615'
Set-TextValue $ws.Range('I182') '1184868481'
Set-TextValue $ws.Range('I199') '1427344720'
Set-TextValue $ws.Range('I200') '1427344720'
Set-TextValue $ws.Range('I205') '1396046694'
Set-TextValue $ws.Range('I214') '1255355277'
Set-TextValue $ws.Range('I215') '1255355277'
Set-TextValue $ws.Range('I220') '1902903966'
Set-TextValue $ws.Range('I225') '1871597435'
Set-TextValue $ws.Range('I235') '1982808473'
Set-TextValue $ws.Range('I238') '1366642415'
Set-TextValue $ws.Range('I239') '1366642415'
Set-TextValue $ws.Range('I242') '1841387941'
Set-TextValue $ws.Range('I248') '1144541905'
$ws.Range('I273').Value = 'UNKNOWN'
Set-TextValue $ws.Range('I280') '1871770123'
Set-TextValue $ws.Range('I287') '1205038296'
Set-TextValue $ws.Range('I302') '1942620232'
Set-TextValue $ws.Range('I303') '1942620232'
Set-TextValue $ws.Range('I341') '1114148855'
Set-TextValue $ws.Range('I342') '1114148855'
Set-TextValue $ws.Range('I358') '1376506337'
Set-TextValue $ws.Range('I359') '1376506337'
Set-TextValue $ws.Range('I396') '1497067235'
Set-TextValue $ws.Range('I397') '1497067235'
Set-TextValue $ws.Range('I406') '1710272521'
Set-TextValue $ws.Range('I407') '1710272521'
Set-TextValue $ws.Range('I438') '1952307233'
$ws.Range('D462').Value = 'Physician'
$ws.Range('E462').Value = 'Physician'
$ws.Range('G462').Value = 'PHY'
Set-TextValue $ws.Range('I462') '1619401239'
$ws.Range('J462').Value = 'Saint Thomas Health Svcs-TN'
$ws.Range('D468').Value = 'Physician'
$ws.Range('E468').Value = 'Physician'
Set-TextValue $ws.Range('I468') '1164426979'
$ws.Range('J468').Value = 'Saint Thomas Health Svcs-TN'
$ws.Range('D469').Value = 'Physician'
$ws.Range('E469').Value = 'Physician'
Set-TextValue $ws.Range('I469') '1164426979'
$ws.Range('J469').Value = 'Saint Thomas Health Svcs-TN'
$ws.Range('D476').Value = 'Non-HR Contractor'
$ws.Range('E476').Value = 'Non-HR Contractor'
$ws.Range('J476').Value = 'Saint Thomas Health Svcs-TN'
$ws.Range('D477').Value = 'Non-HR Contractor'
$ws.Range('E477').Value = 'Non-HR Contractor'
$ws.Range('J477').Value = 'Saint Thomas Health Svcs-TN'
$ws.Range('D482').Value = 'Physician'
$ws.Range('E482').Value = 'Physician'
$ws.Range('G482').Value = 'PHY'
Set-TextValue $ws.Range('I482') '1093912487'
$ws.Range('J482').Value = 'Saint Thomas Health Svcs-TN'
$ws.Range('D483').Value = 'Physician'
$ws.Range('E483').Value = 'Physician'
$ws.Range('G483').Value = 'PHY'
Set-TextValue $ws.Range('I483') '1093912487'
$ws.Range('J483').Value = 'Saint Thomas Health Svcs-TN'
$ws.Range('D487').Value = 'Employee'
$ws.Range('E487').Value = 'Employee'
$ws.Range('G487').Value = 'NUR'
Set-TextValue $ws.Range('I487') '1336411313'
$ws.Range('J487').Value = 'Saint Thomas Health Svcs-TN'
$ws.Range('D488').Value = 'Employee'
$ws.Range('E488').Value = 'Employee'
$ws.Range('G488').Value = 'NUR'
Set-TextValue $ws.Range('I488') '1336411313'
$ws.Range('J488').Value = 'Saint Thomas Health Svcs-TN'
$ws.Range('D492').Value = 'Physician'
$ws.Range('E492').Value = 'Physician'
Set-TextValue $ws.Range('I492') '1649333931'
$ws.Range('J492').Value = 'Saint Thomas Health Svcs-TN'
$ws.Range('D493').Value = 'Physician'
$ws.Range('E493').Value = 'Physician'
Set-TextValue $ws.Range('I493') '1649333931'
$ws.Range('J493').Value = 'Saint Thomas Health Svcs-TN'
$ws.Range('D499').Value = 'Employee'
$ws.Range('E499').Value = 'Employee'
$ws.Range('G499').Value = 'NUR'
Set-TextValue $ws.Range('I499') '1942877618'
$ws.Range('J499').Value = 'Saint Thomas Health Svcs-TN'
$ws.Range('D500').Value = 'Employee'
$ws.Range('E500').Value = 'Employee'
$ws.Range('G500').Value = 'NUR'
Set-TextValue $ws.Range('I500') '1942877618'
$ws.Range('J500').Value = 'Saint Thomas Health Svcs-TN'
